# Updated cryptos list values (prices / 1h-volume / two name-link swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns to stay text so numeric-looking strings
# (e.g. "583.88", "0.0000230") are not auto-converted to numbers.
$ws.Columns.Item(2).NumberFormat = "@"
$ws.Columns.Item(3).NumberFormat = "@"
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

$ws.Range("D2").Value = "62.643.86"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.017.47"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "583.88"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "147.54"
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("D9").Value = "3.015.34"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").Value = "35.02"
$ws.Range("E14").Value = "  -5.08%  "
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "3.527.06"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.654.94"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "7.03"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "3.026.86"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "465.96"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").Value = "13.95"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "0.687"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").Value = "7.35"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "80.49"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "12.37"
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").Value = "10.33"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.63"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "27.62"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "0.0₃0797"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "5.75"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").Value = "2.13"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "50.21"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "8.97"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -10.90%  "
$ws.Range("D42").Value = "423.27"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "0.278"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "2.796.99"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "0.0355"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "37.53"
$ws.Range("E47").Value = "  -8.89%  "
$ws.Range("D48").Value = "128.89"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D50").Value = "24.14"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("E51").Value = "  -0.66%  "
